$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

function Set-TextValue($range, $text) {
    # Force the value to be stored as text even when it looks like a
    # number (e.g. "34562"), without leaving a permanent number-format
    # override on the cell: apply a text format, assign the value, then
    # clear the format override again so the cell's style reverts back
    # to the sheet's default style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = ""
}

Set-TextValue $ws.Range("A10") "zechaad"
Set-TextValue $ws.Range("B10") "1234%z"
Set-TextValue $ws.Range("A11") "xds"
Set-TextValue $ws.Range("B11") "34562"
Set-TextValue $ws.Range("A12") "zechariah"
Set-TextValue $ws.Range("B12") "1234"
